# Update numeric values on the "d2" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("d2")

$ws.Range("A5").Value = 2623.818181818182
$ws.Range("B5").Value = 13089.7

$ws.Range("A7").Value = 636.93181818181802
$ws.Range("B7").Value = 4725.9666666666672

$ws.Range("A25").Value = 278.27272727272725
$ws.Range("B25").Value = 237.54545454545453
$ws.Range("C25").Value = 35.272727272727273

# Make "d2" the active sheet/tab and set its selection
$ws.Activate() | Out-Null
$ws.Range("E4:K21").Select() | Out-Null
